$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the old column B (category), making room for
# the new "designator" and "collection_url" columns.
$ws.Columns("B:C").Insert()

# --- New headers (row 1) ---
$ws.Cells.Item(1,2).Value = "designator"
$ws.Cells.Item(1,3).Value = "collection_url"

# --- Sample data row (row 2), columns B..J ---
$ws.Cells.Item(2,2).Value = "TC"
$ws.Cells.Item(2,3).Value = "thiscollection.org"
$ws.Cells.Item(2,4).Value = "firstcategory"
$ws.Cells.Item(2,5).Value = "cat"
$ws.Cells.Item(2,6).Value = "categorical article"
$ws.Cells.Item(2,7).Value = "categoricalarticle.com"
$ws.Cells.Item(2,8).Value = "Tom Cat"
$ws.Cells.Item(2,9).Value = "Collections"
$ws.Cells.Item(2,10).Value = "2/52/2222"

# --- New trailing "other" column (M) ---
$ws.Cells.Item(1,13).Value = "other"
$ws.Cells.Item(2,13).Value = "accompanyingsong.mp3"

# --- Remaining row 2 values ---
$ws.Cells.Item(2,1).Value = "this collection"
$ws.Cells.Item(2,11).Value = '["cats", "categories", "Tom Cat"]'
$ws.Cells.Item(2,12).Value = "catfile.ris"

# Apply the bold/boxed/centered header formatting (matching the rest of row 1)
# to the new header cells by copying the existing header format, so the same
# style record is reused rather than a near-duplicate being created.
$ws.Range("A1").Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)
$ws.Range("M1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column C (collection_url) width, matching the authored layout.
$ws.Columns("C").ColumnWidth = 11.17

# Selection as left by the author.
[void]$ws.Range("R27").Select()
